# Updated cryptos list (price + volume(1h) refresh, plus a few coins that
# changed their relative rank and swapped rows: Chainlink/Litecoin/TRON in
# rows 14-16, and VeChain/Hedera in rows 38-39).
#
# NOTE: many "Price" values look numeric (e.g. "1.000", "0.4646") but must
# stay as literal text, matching the source data's inlineStr cells. A
# leading apostrophe forces Excel to store them as text (quote-prefixed)
# instead of auto-converting to a number and dropping formatting such as
# trailing/leading zeros. Values that already contain two '.' characters
# (e.g. "26.847.60") are never auto-parsed as numbers, so no apostrophe is
# needed for those.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.847.60"
$ws.Range("E2").Value = "  -1.93%  "
$ws.Range("D3").Value = "1.802.54"
$ws.Range("E3").Value = "  -1.38%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'309.32"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "'0.4646"
$ws.Range("E7").Value = "  +3.77%  "
$ws.Range("D8").Value = "'0.3693"
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("D9").Value = "'0.07365"
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").Value = "'0.8678"
$ws.Range("E10").Value = "  -2.23%  "
$ws.Range("D11").Value = "'20.34"
$ws.Range("E11").Value = "  -3.45%  "
$ws.Range("D12").Value = "1.875.45"
$ws.Range("E12").Value = "  +2.67%  "
$ws.Range("D13").Value = "'5.356"
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D14").Value = "'92.05"
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'6.516"
$ws.Range("E15").Value = "  -3.68%  "
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "'0.07046"
$ws.Range("E16").Value = "  -1.05%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "'0.000008710"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D20").Value = "'14.66"
$ws.Range("E20").Value = "  -3.42%  "
$ws.Range("D21").Value = "26.836.04"
$ws.Range("E21").Value = "  -1.99%  "
$ws.Range("D22").Value = "'5.296"
$ws.Range("E22").Value = "  -2.11%  "
$ws.Range("E23").Value = "  -3.47%  "
$ws.Range("D24").Value = "2.076.65"
$ws.Range("E24").Value = "  +0.90%  "
$ws.Range("D25").Value = "'1.903"
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("D26").Value = "'151.40"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'18.39"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").Value = "'2.132"
$ws.Range("E28").Value = "  -7.88%  "
$ws.Range("D29").Value = "'5.242"
$ws.Range("E29").Value = "  -2.97%  "
$ws.Range("D30").Value = "'115.89"
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").Value = "'0.08913"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "'0.7572"
$ws.Range("E32").Value = "  -4.38%  "
$ws.Range("D33").Value = "'2.934"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").Value = "'1.149"
$ws.Range("E34").Value = "  -4.98%  "
$ws.Range("D35").Value = "'4.448"
$ws.Range("E35").Value = "  -3.39%  "
$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "'1.102"
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.05243"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "'0.01942"
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("E40").Value = "  +2.36%  "
$ws.Range("D41").Value = "'7.202"
$ws.Range("D42").Value = "'0.5269"
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("D43").Value = "'2.347"
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("D44").Value = "'0.1661"
$ws.Range("E44").Value = "  -3.57%  "
$ws.Range("D45").Value = "'8.477"
$ws.Range("E45").Value = "  -2.51%  "
$ws.Range("D46").Value = "'0.5003"
$ws.Range("E46").Value = "  -2.70%  "
$ws.Range("D47").Value = "'10.25"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("D48").Value = "'104.03"
$ws.Range("E48").Value = "  -1.33%  "
$ws.Range("D49").Value = "'0.9998"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "'1.663"
$ws.Range("E50").Value = "  -2.27%  "
$ws.Range("D51").Value = "'0.06283"
$ws.Range("E51").Value = "  -2.00%  "
